$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row: B6 gets a new shared string noting the error-handling catch
# (this is the 5th unique string added to sharedStrings.xml, bumping
# dimension from A3:B5 to A3:B6)
$ws.Range("B6").Value = "Catch de Error en la autenticacion"

# Column B now needs to be wide enough to show that text comfortably -
# mirrors the new <col> entry with a best-fit width.
$ws.Columns.Item(2).ColumnWidth = 38

# Selection moves on, as if the user tabbed past the new row.
$ws.Range("B9").Select()
